# "Assigned achievements to everyone" — add a new column E that assigns
# one of three names (will / zach / steven) to each of the 30 achievement
# rows, cycling through the three names in order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @("will ", "zach ", "steven ")

for ($r = 1; $r -le 30; $r++) {
    $name = $names[($r - 1) % 3]
    $ws.Cells.Item($r, 5).Value = $name
}

# Move / leave the active selection on A5, matching the saved view state.
$ws.Range("A5").Select() | Out-Null
